$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text, preserving numeric-looking strings (e.g. "12", "70.58")
# as literal text instead of letting COM auto-convert them to real numbers.
function Set-TextCell($sheet, $addr, $val) {
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $sheet.Range($addr).Value = "'" + $val
    } else {
        $sheet.Range($addr).Value = $val
    }
}

$newRows = @(
    @{ Row=3; A=' Oct 6 2020'; B=' Abu Dhabi'; C='Mumbai won by 57 runs'; D='Mumbai Indians'; E='Rajasthan Royals'; F='Krunal Pandya '; G='12'; H='17'; I='0'; J='1'; K='70.58' }
    @{ Row=4; A=' Oct 18 2020'; B=' Dubai (DSC)'; C='Match tied (Kings XI won the one-over eliminator)'; D='Mumbai Indians'; E='Kings XI Punjab'; F='Krunal Pandya '; G='34'; H='30'; I='4'; J='1'; K='113.33' }
    @{ Row=5; A=' Oct 25 2020'; B=' Abu Dhabi'; C='Royals won by 8 wickets (with 10 balls remaining)'; D='Mumbai Indians'; E='Rajasthan Royals'; F='Krunal Pandya '; G='3'; H='4'; I='0'; J='0'; K='75.00' }
    @{ Row=6; A=' Oct 4 2020'; B=' Sharjah'; C='Mumbai won by 34 runs'; D='Mumbai Indians'; E='Sunrisers Hyderabad'; F='Krunal Pandya '; G='20'; H='4'; I='2'; J='2'; K='500.00' }
    @{ Row=7; A=' Oct 11 2020'; B=' Abu Dhabi'; C='Mumbai won by 5 wickets (with 2 balls remaining)'; D='Mumbai Indians'; E='Delhi Capitals'; F='Krunal Pandya '; G='12'; H='7'; I='2'; J='0'; K='171.42' }
    @{ Row=8; A=' Nov 5 2020'; B=' Dubai (DSC)'; C='Mumbai won by 57 runs'; D='Mumbai Indians'; E='Delhi Capitals'; F='Krunal Pandya '; G='13'; H='10'; I='0'; J='1'; K='130.00' }
    @{ Row=9; A=' Sep 19 2020'; B=' Abu Dhabi'; C='Super Kings won by 5 wickets (with 4 balls remaining)'; D='Mumbai Indians'; E='Chennai Super Kings'; F='Krunal Pandya '; G='3'; H='3'; I='0'; J='0'; K='100.00' }
    @{ Row=10; A=' Nov 3 2020'; B=' Sharjah'; C='Sunrisers won by 10 wickets (with 17 balls remaining)'; D='Mumbai Indians'; E='Sunrisers Hyderabad'; F='Krunal Pandya '; G='0'; H='3'; I='0'; J='0'; K='0.00' }
    @{ Row=11; A=' Sep 28 2020'; B=' Dubai (DSC)'; C='Match tied (RCB won the one-over eliminator)'; D='Mumbai Indians'; E='Royal Challengers Bangalore'; F='Krunal Pandya '; G='0'; H='0'; I='0'; J='0'; K='-' }
    @{ Row=12; A=' Sep 23 2020'; B=' Abu Dhabi'; C='Mumbai won by 49 runs'; D='Mumbai Indians'; E='Kolkata Knight Riders'; F='Krunal Pandya '; G='1'; H='3'; I='0'; J='0'; K='33.33' }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $addrA = "A" + $rowNum
    Set-TextCell $ws $addrA $r.A
    $addrB = "B" + $rowNum
    Set-TextCell $ws $addrB $r.B
    $addrC = "C" + $rowNum
    Set-TextCell $ws $addrC $r.C
    $addrD = "D" + $rowNum
    Set-TextCell $ws $addrD $r.D
    $addrE = "E" + $rowNum
    Set-TextCell $ws $addrE $r.E
    $addrF = "F" + $rowNum
    Set-TextCell $ws $addrF $r.F
    $addrG = "G" + $rowNum
    Set-TextCell $ws $addrG $r.G
    $addrH = "H" + $rowNum
    Set-TextCell $ws $addrH $r.H
    $addrI = "I" + $rowNum
    Set-TextCell $ws $addrI $r.I
    $addrJ = "J" + $rowNum
    Set-TextCell $ws $addrJ $r.J
    $addrK = "K" + $rowNum
    Set-TextCell $ws $addrK $r.K
}
